$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) mirroring the existing year-header / data-row
# formatting already used for column R.

# Header cell S4: same format as R4 (year header style), value 2022
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Data cell S5: same base format as P5/Q5/R5 (data style), but with a
# "0.0" number format, value 42
$ws.Range("P5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").NumberFormat = "0.0"
$ws.Range("S5").Value = 42

# Move the active selection to U4, matching the author's final cursor spot
$ws.Range("U4").Select() | Out-Null
